$d = $word.ActiveDocument

# Locate the end of "...exe file." inside the "Alternative way to obtain a
# root shell" paragraph so we can split the run there and insert a new
# sentence before " To do it, I installed on my Windows 11 machine the ".
$find = $d.Content.Find
$find.Text = ", I can convert that python file in an exe file."
$found = $find.Execute()

if ($found) {
    $splitPoint = $d.Range($find.Parent.End, $find.Parent.End)

    # Insert the new sentence right after "exe file." and before " To do it...".
    $splitPoint.InsertAfter(" Actually, it is better to generate a python shell using msfvenom and convert it, so you can control the payload to use.")

    # Toggling formatting and reverting it forces the engine to keep this
    # insertion as its own run instead of silently re-merging it with the
    # (identically formatted) neighboring runs.
    $splitPoint.Font.Bold = 1
    $splitPoint.Font.Bold = 0
}
